$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A1:A23 with values 0..22
for ($i = 0; $i -le 22; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $i
}

# Update the active selection to E16
$ws.Range("E16").Select() | Out-Null
